$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/week number & report date range) ---
$ws.Range("A8").Characters(21,2).Text = "45"
$ws.Range("C9").Value = "Report Covering the Week  11/7/2022  Through  11/13/2022"

# --- Crime statistics table updates (rows 14-29) ---
$ws.Range("C14").NumberFormat = "#,##0"
$ws.Range("C14").Value = 1
$ws.Range("F14").NumberFormat = "#,##0"
$ws.Range("F14").Value = 1
$ws.Range("I14").Value = 5
$ws.Range("K14").Value = 400
$ws.Range("L14").Value = 66.666666666666
$ws.Range("M14").Value = 400
$ws.Range("N14").Value = -37.5
$ws.Range("F15").Value = 3
$ws.Range("H15").Value = 200
$ws.Range("I15").Value = 10
$ws.Range("K15").Value = -9.090909090909
$ws.Range("L15").Value = 66.666666666666
$ws.Range("M15").Value = 25
$ws.Range("N15").Value = -56.521739130434
$ws.Range("C16").Value = 8
$ws.Range("E16").Value = 166.666666666667
$ws.Range("F16").Value = 15
$ws.Range("G16").Value = 11
$ws.Range("H16").Value = 36.363636363636
$ws.Range("I16").Value = 138
$ws.Range("J16").Value = 116
$ws.Range("K16").Value = 18.965517241379
$ws.Range("L16").Value = 26.605504587156
$ws.Range("M16").Value = 46.808510638297
$ws.Range("N16").Value = -74.349442379182
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 200
$ws.Range("F17").Value = 16
$ws.Range("G17").Value = 16
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 121
$ws.Range("J17").Value = 117
$ws.Range("K17").Value = 3.418803418803
$ws.Range("L17").Value = 28.723404255319
$ws.Range("M17").Value = 5.217391304347
$ws.Range("N17").Value = -45.495495495495
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 50
$ws.Range("F18").Value = 14
$ws.Range("G18").Value = 12
$ws.Range("H18").Value = 16.666666666666
$ws.Range("I18").Value = 150
$ws.Range("J18").Value = 113
$ws.Range("K18").Value = 32.743362831858
$ws.Range("L18").Value = 0.671140939597
$ws.Range("M18").Value = 68.539325842696
$ws.Range("N18").Value = -72.776769509981
$ws.Range("C19").Value = 12
$ws.Range("D19").Value = 16
$ws.Range("E19").Value = -25
$ws.Range("F19").Value = 71
$ws.Range("G19").Value = 63
$ws.Range("H19").Value = 12.698412698412
$ws.Range("I19").Value = 624
$ws.Range("J19").Value = 453
$ws.Range("K19").Value = 37.748344370860
$ws.Range("L19").Value = 34.193548387096
$ws.Range("M19").Value = 10.442477876106
$ws.Range("N19").Value = -15.217391304347
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 10
$ws.Range("G20").Value = 10
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 73
$ws.Range("J20").Value = 54
$ws.Range("K20").Value = 35.185185185185
$ws.Range("L20").Value = 73.809523809523
$ws.Range("M20").Value = 87.179487179487
$ws.Range("N20").Value = -83.885209713024
$ws.Range("C21").Value = 35
$ws.Range("D21").Value = 27
$ws.Range("E21").Value = 29.629629629629
$ws.Range("F21").Value = 130
$ws.Range("G21").Value = 113
$ws.Range("H21").Value = 15.044247787610
$ws.Range("I21").Value = 1121
$ws.Range("J21").Value = 865
$ws.Range("K21").Value = 29.595375722543
$ws.Range("L21").Value = 29.147465437788
$ws.Range("M21").Value = 23.051591657519
$ws.Range("N21").Value = -55.709205847491
$ws.Range("F22").Value = 5
$ws.Range("H22").Value = 25
$ws.Range("I22").Value = 29
$ws.Range("J22").Value = 18
$ws.Range("K22").Value = 61.111111111111
$ws.Range("L22").Value = 141.666666666667
$ws.Range("M22").Value = 123.076923076923
$ws.Range("D23").Value = 1
$ws.Range("F23").Value = 2
$ws.Range("H23").Value = -66.666666666666
$ws.Range("J23").Value = 59
$ws.Range("K23").Value = -18.644067796610
$ws.Range("L23").Value = 9.090909090909
$ws.Range("C24").Value = 12
$ws.Range("D24").Value = 22
$ws.Range("E24").Value = -45.454545454545
$ws.Range("F24").Value = 59
$ws.Range("G24").Value = 59
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 653
$ws.Range("J24").Value = 495
$ws.Range("K24").Value = 31.919191919191
$ws.Range("L24").Value = 1.397515527950
$ws.Range("M24").Value = -13.049267643142
$ws.Range("C25").Value = 4
$ws.Range("D25").Value = 10
$ws.Range("E25").Value = -60
$ws.Range("F25").Value = 27
$ws.Range("G25").Value = 31
$ws.Range("H25").Value = -12.903225806451
$ws.Range("I25").Value = 334
$ws.Range("J25").Value = 285
$ws.Range("K25").Value = 17.192982456140
$ws.Range("L25").Value = 53.211009174311
$ws.Range("M25").Value = 15.972222222222
$ws.Range("D26").NumberFormat = "#,##0"
$ws.Range("D26").Value = 1
$ws.Range("E26").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 4
$ws.Range("H26").Value = 100
$ws.Range("I26").Value = 16
$ws.Range("J26").Value = 14
$ws.Range("K26").Value = 14.285714285714
$ws.Range("L26").Value = 100
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 0
$ws.Range("I27").Value = 51
$ws.Range("J27").Value = 46
$ws.Range("K27").Value = 10.869565217391
$ws.Range("L27").Value = 37.837837837837
$ws.Range("C28").NumberFormat = "#,##0"
$ws.Range("C28").Value = 3
$ws.Range("F28").NumberFormat = "#,##0"
$ws.Range("F28").Value = 3
$ws.Range("I28").Value = 6
$ws.Range("K28").Value = 200
$ws.Range("L28").Value = 50
$ws.Range("M28").Value = 20
$ws.Range("N28").Value = -40
$ws.Range("C29").NumberFormat = "#,##0"
$ws.Range("C29").Value = 1
$ws.Range("F29").NumberFormat = "#,##0"
$ws.Range("F29").Value = 1
$ws.Range("I29").Value = 3
$ws.Range("K29").Value = 50
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = 0
$ws.Range("N29").Value = -70
